# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the refreshed counts from the latest data pull.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value for column F
$exhibitUpdates = @{
    2  = 1135
    4  = 256
    6  = 12121
    7  = 50
    9  = 11892
    10 = 4778
    11 = 592
    13 = 30
    15 = 88
    16 = 936
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$allUpdates = @{
    2  = 1135
    4  = 256
    8  = 12121
    9  = 50
    11 = 11892
    12 = 4778
    13 = 592
    15 = 30
    17 = 88
    18 = 936
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
